$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = "TV da sala"
$ws.Range("B1").Value = "Televisor"
$ws.Range("D1").Value = 0
$ws.Range("E1").Value = $false

$ws.Range("A2").Value = "Tv da cozinha"
$ws.Range("B2").Value = "Televisor"
$ws.Range("D2").Value = 0
$ws.Range("E2").Value = $false

$ws.Range("C1").ClearContents()
